# Add new Block property rows (11-22) with their Id / Type / flags, mirroring
# the existing rows (2-10) in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newIds = @(
    "Grass1",
    "Grass2",
    "Grass3",
    "Grass4",
    "Grass5",
    "Crack1",
    "Crack2",
    "Crack3",
    "Crack4",
    "Crack5",
    "Treasure1 ",
    "Treasure2"
)

$startRow = 11
for ($i = 0; $i -lt $newIds.Count; $i++) {
    $r = $startRow + $i

    $ws.Range("A$r").Value = $newIds[$i]
    # Rows 21 and 22 keep the default (general) cell format, matching the
    # original author's data entry; every other row's Id column is
    # formatted as text like the pre-existing rows.
    if ($r -ne 21 -and $r -ne 22) {
        $ws.Range("A$r").NumberFormat = "@"
    }

    $ws.Range("B$r").Value = "string"
    $ws.Range("B$r").NumberFormat = "@"

    $ws.Range("C$r").Value = $false
    $ws.Range("D$r").Value = $false
    $ws.Range("E$r").Value = $false
    $ws.Range("F$r").Value = $true

    $ws.Range("G$r").Value = 0
    $ws.Range("H$r").Value = 0

    $ws.Range("I$r").Value = "Friend"
    $ws.Range("I$r").NumberFormat = "@"
}

$ws.Range("H24").Select()
